$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log row (time, cost) after the existing data (rows 1-26).
$ws.Range("A27").Value = "2023-12-06 15:13:33"
$ws.Range("B27").Value = 0.0002
